# "4.0.3 model and data"
#
# The CID "Check Input Data" workbook's InputData file-list is updated:
# the transportation vehicle-type-quantity files `trans/BVTQaZ/BVTQaZ.csv`
# and `trans/VTQaZ/VTQaZ.csv` (each listed once on the "Boolean" sheet)
# are each split into six per-vehicle-class CSV files (LDVs, HDVs,
# aircraft, rail, ships, motorbikes).

$wb = $excel.ActiveWorkbook

$wsAbout      = $wb.Worksheets.Item("About")
$wsInteger    = $wb.Worksheets.Item("Integer")
$wsBoolean    = $wb.Worksheets.Item("Boolean")
$wsSubscript  = $wb.Worksheets.Item("Subscript")

# ---------------------------------------------------------------------
# "Boolean" sheet: split the two combined transportation csv rows into
# six rows apiece (one per vehicle class).
# ---------------------------------------------------------------------

# Before the edit, row 17 holds trans/BVTQaZ/BVTQaZ.csv. Insert five more
# rows right after it so the single row can become six.
$wsBoolean.Range("A18:A22").EntireRow.Insert()

$wsBoolean.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBoolean.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBoolean.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBoolean.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBoolean.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBoolean.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# Rows now read (23 BVTStL, 24 PVTStL, 25 SRPbVT, 26 VTQaZ, 27 VTStFES).
# Row 26 holds trans/VTQaZ/VTQaZ.csv; insert five more rows after it.
$wsBoolean.Range("A27:A31").EntireRow.Insert()

$wsBoolean.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBoolean.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBoolean.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBoolean.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBoolean.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBoolean.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# A handful of blank, formatted rows trail the list.
$wsBoolean.Range("A33:A38").Font.Name = "Calibri"
$wsBoolean.Range("A33:A38").Font.Size = 11

# ---------------------------------------------------------------------
# View / selection state: "About" becomes the active tab; "Integer" and
# "Boolean" keep their own remembered selections; "Subscript" untouched.
# ---------------------------------------------------------------------

[void]$wsInteger.Select()
[void]$wsInteger.Range("A13").Select()

[void]$wsBoolean.Select()
[void]$wsBoolean.Range("A32").Select()

[void]$wsAbout.Select()
